$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.690.45"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.256.76"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'606.20"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'158.87"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.255.49"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "'5.95"
$ws.Range("E11").Value = "  +4.55%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "'39.53"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "3.790.90"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "66.711.22"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'7.39"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.255.88"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "'508.31"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "'15.40"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'0.753"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").Value = "'8.08"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "'14.82"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'86.48"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "'0.156"
$ws.Range("E26").Value = "  +74.27%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "'2.39"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "'6.88"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'2.86"
$ws.Range("E32").Value = "  -7.15%  "
$ws.Range("D33").Value = "'28.23"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("D36").Value = "'6.44"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "0.0₃0801"
$ws.Range("E37").Value = "  +16.74%  "
$ws.Range("E38").Value = "  +19.91%  "
$ws.Range("D39").Value = "'55.52"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'495.50"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'0.0428"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "'0.128"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'8.85"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "'0.296"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "2.947.70"
$ws.Range("D47").Value = "'28.66"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'121.63"
$ws.Range("E51").Value = "  -0.77%  "
